$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 2527
$ws.Range("F3").Value = 547
$ws.Range("F5").Value = 282
$ws.Range("F6").Value = 179
$ws.Range("F7").Value = 448
$ws.Range("F8").Value = 1177
$ws.Range("F9").Value = 529
$ws.Range("F10").Value = 287
$ws.Range("F11").Value = 110
$ws.Range("F12").Value = 341
$ws.Range("F13").Value = 5454
$ws.Range("F14").Value = 44
$ws.Range("F15").Value = 1634
$ws.Range("F16").Value = 3965
$ws.Range("F17").Value = 396
$ws.Range("F20").Value = 4504
$ws.Range("F21").Value = 5896
$ws.Range("F23").Value = 1015
$ws.Range("F24").Value = 646
$ws.Range("F25").Value = 3648
$ws.Range("F26").Value = 460
$ws.Range("F28").Value = 179
$ws.Range("F30").Value = 954
$ws.Range("F31").Value = 1348
$ws.Range("F32").Value = 440
$ws.Range("F33").Value = 502
$ws.Range("F34").Value = 1544
$ws.Range("F35").Value = 184
$ws.Range("F36").Value = 1619
$ws.Range("F37").Value = 149
$ws.Range("F38").Value = 1
$ws.Range("F39").Value = 1063
$ws.Range("F41").Value = 1342
$ws.Range("F42").Value = 598
$ws.Range("F43").Value = 85
$ws.Range("F44").Value = 196
$ws.Range("F45").Value = 2738
$ws.Range("F46").Value = 117
$ws.Range("F47").Value = 252
$ws.Range("F48").Value = 394
$ws.Range("F49").Value = 3849

$ws = $wb.Worksheets.Item("演出")
$ws.Range("G2").Value = 180
$ws.Range("F5").Value = 1160
$ws.Range("F6").Value = 39
$ws.Range("F14").Value = 11
$ws.Range("F18").Value = 7

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 3633

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 2527
$ws.Range("F4").Value = 547
$ws.Range("F6").Value = 282
$ws.Range("F7").Value = 1160
$ws.Range("F8").Value = 179
$ws.Range("F9").Value = 448
$ws.Range("F10").Value = 1177
$ws.Range("F11").Value = 529
$ws.Range("F12").Value = 287
$ws.Range("F13").Value = 110
$ws.Range("F14").Value = 341
$ws.Range("F15").Value = 5454
$ws.Range("F17").Value = 1634
$ws.Range("F18").Value = 4504
$ws.Range("F19").Value = 5896
$ws.Range("F21").Value = 1015
$ws.Range("F22").Value = 646
$ws.Range("F23").Value = 3648
$ws.Range("F24").Value = 460
$ws.Range("F26").Value = 179
$ws.Range("F28").Value = 1348
$ws.Range("F29").Value = 440
$ws.Range("F30").Value = 502
$ws.Range("F32").Value = 1544
$ws.Range("F33").Value = 184
$ws.Range("F34").Value = 1619
$ws.Range("F36").Value = 1063
$ws.Range("F37").Value = 7
$ws.Range("F38").Value = 598
$ws.Range("F41").Value = 85
$ws.Range("F43").Value = 2738
$ws.Range("F45").Value = 117
$ws.Range("F46").Value = 252
$ws.Range("F47").Value = 394
$ws.Range("F49").Value = 3849
